$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column C. This pushes the existing headers
# (old C:AB = "Diad2_HB2_prom_ratio" .. "Mean_HB_prom") one column to the
# right (new D:AC), and leaves a blank column C ready for the new header.
$ws.Columns("C:C").Insert()

# Give the new column C the same header look/format as the rest of the
# header row (bold font, borders, centered/top aligned) by copying B1's
# formatting, then set its text.
$ws.Range("B1").Copy() | Out-Null
$ws.Range("C1").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$ws.Range("C1").Value = "rays_present"

# Append two brand-new header columns at the end of the row.
$ws.Range("AD1").Value = "Diad2_prom/std_betweendiads"
$ws.Range("AE1").Value = "Left_vs_Right"
$ws.Range("B1").Copy() | Out-Null
$ws.Range("AD1:AE1").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

$excel.CutCopyMode = 0
